$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2673
$ws.Range("B3").Value = 2408
$ws.Range("B4").Value = 3315
$ws.Range("B6").Value = 628
$ws.Range("B8").Value = 66872
$ws.Range("B9").Value = 4209
$ws.Range("B11").Value = 10351
$ws.Range("B12").Value = 4650
$ws.Range("B13").Value = 212
$ws.Range("B14").Value = 674
$ws.Range("B15").Value = 11833
$ws.Range("B17").Value = 2612
$ws.Range("B18").Value = 24483
$ws.Range("B22").Value = 13151
$ws.Range("B23").Value = 8790
$ws.Range("B24").Value = 734
$ws.Range("B25").Value = 419114
$ws.Range("B27").Value = 16886
$ws.Range("B28").Value = 162
$ws.Range("B31").Value = 231
$ws.Range("B32").Value = 114
$ws.Range("B33").Value = 1152
$ws.Range("B34").Value = 24523
$ws.Range("B35").Value = 93
$ws.Range("B37").Value = 27004
$ws.Range("B39").Value = 76867
$ws.Range("B41").Value = 148
$ws.Range("B42").Value = 772
$ws.Range("B43").Value = 3365
$ws.Range("B45").Value = 7388
$ws.Range("B46").Value = 713
$ws.Range("B47").Value = 332
$ws.Range("B48").Value = 29608
$ws.Range("B49").Value = 2496
$ws.Range("B50").Value = 148
$ws.Range("B52").Value = 3517
$ws.Range("B53").Value = 19061
$ws.Range("B54").Value = 13779
$ws.Range("B55").Value = 2146
$ws.Range("B58").Value = 1193
$ws.Range("B60").Value = 3840
$ws.Range("B61").Value = 3
$ws.Range("B62").Value = 922
$ws.Range("B63").Value = 106262
$ws.Range("B64").Value = 142
$ws.Range("B66").Value = 4245
$ws.Range("B67").Value = 84659
$ws.Range("B69").Value = 10910
$ws.Range("B71").Value = 7695
$ws.Range("B72").Value = 149
$ws.Range("B74").Value = 314
$ws.Range("B75").Value = 263
$ws.Range("B76").Value = 5585
$ws.Range("B77").Value = 28403
$ws.Range("B79").Value = 238270
$ws.Range("B80").Value = 46663
$ws.Range("B81").Value = 74241
$ws.Range("B82").Value = 15702
$ws.Range("B83").Value = 4918
$ws.Range("B84").Value = 6375
$ws.Range("B85").Value = 122470
$ws.Range("B86").Value = 801
$ws.Range("B87").Value = 10712
$ws.Range("B88").Value = 9047
$ws.Range("B89").Value = 3369
$ws.Range("B90").Value = 2865
$ws.Range("B91").Value = 1865
$ws.Range("B92").Value = 2207
$ws.Range("B93").Value = 1628
$ws.Range("B94").Value = 1643
$ws.Range("B96").Value = 2197
$ws.Range("B97").Value = 7436
$ws.Range("B100").Value = 3063
$ws.Range("B101").Value = 58
$ws.Range("B102").Value = 4019
$ws.Range("B103").Value = 801
$ws.Range("B104").Value = 701
$ws.Range("B105").Value = 1152
$ws.Range("B106").Value = 1632
$ws.Range("B107").Value = 82
$ws.Range("B108").Value = 497
$ws.Range("B109").Value = 417
$ws.Range("B113").Value = 218657
$ws.Range("B114").Value = 5929
$ws.Range("B116").Value = 160
$ws.Range("B117").Value = 1533
$ws.Range("B118").Value = 9057
$ws.Range("B119").Value = 820
$ws.Range("B120").Value = 682
$ws.Range("B121").Value = 3579
$ws.Range("B122").Value = 17547
$ws.Range("B128").Value = 2083
$ws.Range("B129").Value = 18797
$ws.Range("B130").Value = 6258
$ws.Range("B132").Value = 6974
$ws.Range("B133").Value = 62976
$ws.Range("B134").Value = 18099
$ws.Range("B135").Value = 69445
$ws.Range("B136").Value = 16989
$ws.Range("B137").Value = 496
$ws.Range("B138").Value = 28799
$ws.Range("B139").Value = 110735
$ws.Range("B142").Value = 75
$ws.Range("B143").Value = 12
$ws.Range("B147").Value = 7045
$ws.Range("B148").Value = 1116
$ws.Range("B149").Value = 6519
$ws.Range("B153").Value = 11946
$ws.Range("B154").Value = 4285
$ws.Range("B156").Value = 745
$ws.Range("B157").Value = 54687
$ws.Range("B159").Value = 78792
$ws.Range("B160").Value = 764
$ws.Range("B162").Value = 214
$ws.Range("B163").Value = 14173
$ws.Range("B164").Value = 10702
$ws.Range("B165").Value = 1639
$ws.Range("B169").Value = 363
$ws.Range("B171").Value = 124
$ws.Range("B172").Value = 196
$ws.Range("B173").Value = 11277
$ws.Range("B174").Value = 42465
$ws.Range("B175").Value = 580901
$ws.Range("B176").Value = 346
$ws.Range("B177").Value = 47717
$ws.Range("B178").Value = 1607
$ws.Range("B179").Value = 127858
$ws.Range("B180").Value = 3032
$ws.Range("B181").Value = 660
$ws.Range("B183").Value = 2263
$ws.Range("B185").Value = 3338
$ws.Range("B186").Value = 1269
$ws.Range("B187").Value = 1256
$ws.Range("B188").Value = 1576
